$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "93.696.82"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.093.47"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "236.95"
$ws.Range("E5").Value = "  -3.37%  "
$ws.Range("D6").Value = "611.56"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("D7").Value = "1.12"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("D8").Value = "0.385"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "0.818"
$ws.Range("E10").Value = "  +9.55%  "
$ws.Range("D11").Value = "3.090.72"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "0.197"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("E13").Value = "  -4.41%  "
$ws.Range("D14").Value = "93.331.53"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "34.51"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").Value = "3.665.76"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "3.101.43"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "3.65"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").Value = "5.89"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "441.26"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("E24").Value = "  -5.59%  "
$ws.Range("D25").Value = "8.26"
$ws.Range("E25").Value = "  +4.83%  "
$ws.Range("D26").Value = "5.62"
$ws.Range("D27").Value = "12.22"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "85.81"
$ws.Range("E28").Value = "  -3.21%  "
$ws.Range("D29").Value = "3.256.80"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "0.243"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "0.180"
$ws.Range("E32").Value = "  +6.45%  "
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  -11.96%  "
$ws.Range("D34").Value = "9.11"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "7.90"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").Value = "  -11.29%  "
$ws.Range("D38").Value = "25.74"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("B39").Value = "MantraDAO"
$ws.Range("C39").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D39").Value = "3.88"
$ws.Range("E39").Value = "  -10.67%  "
$ws.Range("B40").Value = "PancakeSwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D40").Value = "1.88"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").Value = "0.445"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "23.97"
$ws.Range("E42").Value = "  +7.85%  "
$ws.Range("D43").Value = "474.05"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "3.21"
$ws.Range("E46").Value = "  -8.01%  "
$ws.Range("D47").Value = "159.43"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "0.683"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E49").Value = "  -5.80%  "
$ws.Range("D50").Value = "43.77"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  -5.59%  "
